$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4,  17, 3, 11, 1000,  "2025-03-27 08:25:11", 0),
    @(5,  17, 3, 11, 500,   "2025-03-27 08:25:48", 0),
    @(6,  17, 3, 11, 2000,  "2025-03-27 08:32:13", 0),
    @(7,  17, 3, 11, 1000,  "2025-03-27 08:54:02", 0),
    @(8,  17, 3, 11, 5000,  "2025-03-27 08:54:32", 0),
    @(9,  17, 4, 10, 10000, "2025-03-27 15:31:39", 0),
    @(10, 17, 3, 11, 1000,  "2025-03-27 15:32:22", 0),
    @(11, 18, 5, 11, 350,   "2025-03-27 17:21:11", 0)
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
}
